$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: insert $text right before the end of the paragraph currently
# at index $paraIndex (i.e. after any text already there, but before the
# paragraph mark), as a *separate* run. Toggling Bold on/off right after
# the insert keeps the new text from being silently re-merged into the
# neighbouring run during serialization, while leaving the final
# run-properties untouched (Bold is left exactly as it was: off).
# Returns the (start, length) of the inserted text.
# ---------------------------------------------------------------------
function Insert-SplitRun($paraIndex, $text) {
    $r = $d.Paragraphs.Item($paraIndex).Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Collapse(0)
    $start = $r.Start
    $r.InsertAfter($text)
    $len = $text.Length
    $newRng = $d.Range($start, $start + $len)
    $newRng.Font.Bold = 1
    $newRng.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) Drop the _GoBack bookmark that used to sit at the end of the
#    "Update version-name and version-code in Android Manifest" item.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Turn "Update the change notes" into four runs:
#    "Update the " / "document "Release" / " notes" / """
# ---------------------------------------------------------------------
$notesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "Update the change notes") {
        $notesPara = $i
        break
    }
}
if (-not $notesPara) {
    throw "Could not locate the 'Update the change notes' paragraph"
}

$notesRange = $d.Paragraphs.Item($notesPara).Range
$notesRange.MoveEnd(1, -1) | Out-Null
$notesRange.Text = ""

Insert-SplitRun $notesPara "Update the "
Insert-SplitRun $notesPara "document “Release"
Insert-SplitRun $notesPara " notes"
Insert-SplitRun $notesPara "”"

# Make sure every run we just created carries the en-US language tag
# (matches the rest of the document's runs).
$fixRange = $d.Paragraphs.Item($notesPara).Range
$fixRange.MoveEnd(1, -1) | Out-Null
$fixRange.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 3) Add a brand-new list paragraph right after it:
#    "Update the document "" / "Known Bugs and Limitations" / """
# ---------------------------------------------------------------------
$notesRangeFull = $d.Paragraphs.Item($notesPara).Range
$notesRangeFull.InsertParagraphAfter()

$bugsPara = $notesPara + 1
$bugsRange = $d.Paragraphs.Item($bugsPara).Range
$bugsRange.MoveEnd(1, -1) | Out-Null
$bugsRange.Text = ""

Insert-SplitRun $bugsPara "Update the document “"
Insert-SplitRun $bugsPara "Known Bugs and Limitations"
Insert-SplitRun $bugsPara "”"

$fixRange2 = $d.Paragraphs.Item($bugsPara).Range
$fixRange2.MoveEnd(1, -1) | Out-Null
$fixRange2.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 4) Re-create the _GoBack bookmark at the end of the new paragraph.
#    (Document.Bookmarks.Add needs a non-empty Range to anchor
#    correctly in this host, so we wrap a throw-away placeholder
#    character and then delete it -- the bookmark collapses to that
#    exact spot and stays put, same as real Word does when you delete
#    bookmarked text.)
# ---------------------------------------------------------------------
$bmHost = $d.Paragraphs.Item($bugsPara).Range
$bmHost.MoveEnd(1, -1) | Out-Null
$bmHost.Collapse(0)
$phStart = $bmHost.Start
$bmHost.InsertAfter("X")
$phRange = $d.Range($phStart, $phStart + 1)
$d.Bookmarks.Add("_GoBack", $phRange)
$d.Range($phStart, $phStart + 1).Text = ""

Write-Output "done"
